$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H19").Value = 2771.1
$ws_ALC.Range("I19").Value = 5572.75
$ws_ALC.Range("J19").Value = 903.3333
$ws_ALC.Range("K19").Value = 5572.75
$ws_ALC.Range("L19").Value = 903.3333
$ws_ALC.Range("M19").Value = -5397.75
$ws_ALC.Range("N19").Value = -1253.3333
$ws_ALC.Range("H32").Value = 1000
$ws_ALC.Range("J32").Value = 1000
$ws_ALC.Range("L32").Value = 1000
$ws_ALC.Range("N32").Value = -1652
$ws_ALC.Range("H40").Value = 1241.1818
$ws_ALC.Range("I40").Value = 1032
$ws_ALC.Range("K40").Value = 1032
$ws_ALC.Range("M40").Value = -857
$ws_ALC.Range("H51").Value = 16240
$ws_ALC.Range("I51").Value = 20000
$ws_ALC.Range("J51").Value = 12480
$ws_ALC.Range("K51").Value = 20000
$ws_ALC.Range("L51").Value = 12480
$ws_ALC.Range("M51").Value = -19516
$ws_ALC.Range("N51").Value = -13448
$ws_ALC.Range("H53").Value = 2970.818
$ws_ALC.Range("I53").Value = 222.375
$ws_ALC.Range("J53").Value = 10300
$ws_ALC.Range("K53").Value = 222.375
$ws_ALC.Range("L53").Value = 10300
$ws_ALC.Range("M53").Value = 414.625
$ws_ALC.Range("N53").Value = -11574
$ws_ALC.Range("H64").Value = 4142.143
$ws_ALC.Range("I64").Value = 3943.4443
$ws_ALC.Range("J64").Value = 4499.8
$ws_ALC.Range("K64").Value = 3943.4443
$ws_ALC.Range("L64").Value = 4499.8
$ws_ALC.Range("M64").Value = -3695.4443
$ws_ALC.Range("N64").Value = -4995.8
$ws_ALC.Range("H67").Value = 4142.143
$ws_ALC.Range("I67").Value = 3943.4443
$ws_ALC.Range("J67").Value = 4499.8
$ws_ALC.Range("K67").Value = 3943.4443
$ws_ALC.Range("L67").Value = 4499.8
$ws_ALC.Range("M67").Value = -3085.4443
$ws_ALC.Range("N67").Value = -6215.8
$ws_ALC.Range("H74").Value = 17861414
$ws_ALC.Range("I74").Value = 3951.5
$ws_ALC.Range("J74").Value = 25004400
$ws_ALC.Range("K74").Value = 3951.5
$ws_ALC.Range("L74").Value = 25004400
$ws_ALC.Range("M74").Value = -3015.5
$ws_ALC.Range("N74").Value = -25006272
$ws_ALC.Range("H77").Value = 17861414
$ws_ALC.Range("I77").Value = 3951.5
$ws_ALC.Range("J77").Value = 25004400
$ws_ALC.Range("K77").Value = 19757.5
$ws_ALC.Range("L77").Value = 125022000
$ws_ALC.Range("M77").Value = -15077.5
$ws_ALC.Range("N77").Value = -125031360
$ws_ALC.Range("H112").Value = 4445425
$ws_ALC.Range("I112").Value = 375
$ws_ALC.Range("J112").Value = 4831951
$ws_ALC.Range("K112").Value = 1125
$ws_ALC.Range("L112").Value = 14495853
$ws_ALC.Range("M112").Value = -17
$ws_ALC.Range("N112").Value = -14498069
$ws_ARM.Range("H63").Value = 2405852.8
$ws_ARM.Range("I63").Value = 2242.7778
$ws_ARM.Range("J63").Value = 7813975
$ws_ARM.Range("K63").Value = 2242.7778
$ws_ARM.Range("L63").Value = 7813975
$ws_ARM.Range("M63").Value = -1556.7778
$ws_ARM.Range("N63").Value = -7815347
$ws_ARM.Range("H66").Value = 2405852.8
$ws_ARM.Range("I66").Value = 2242.7778
$ws_ARM.Range("J66").Value = 7813975
$ws_ARM.Range("K66").Value = 11213.889
$ws_ARM.Range("L66").Value = 39069875
$ws_ARM.Range("M66").Value = -7781.888999999999
$ws_ARM.Range("N66").Value = -39076739
$ws_ARM.Range("H134").Value = 46063.6
$ws_ARM.Range("J134").Value = 46063.6
$ws_ARM.Range("L134").Value = 46063.6
$ws_ARM.Range("N134").Value = -56203.6
$ws_BSM.Range("H86").Value = 2033.0869
$ws_BSM.Range("I86").Value = 1844.7646
$ws_BSM.Range("J86").Value = 2566.6667
$ws_BSM.Range("K86").Value = 1844.7646
$ws_BSM.Range("L86").Value = 2566.6667
$ws_BSM.Range("M86").Value = -721.7646
$ws_BSM.Range("N86").Value = -4812.6667
$ws_BSM.Range("H89").Value = 2033.0869
$ws_BSM.Range("I89").Value = 1844.7646
$ws_BSM.Range("J89").Value = 2566.6667
$ws_BSM.Range("K89").Value = 9223.823
$ws_BSM.Range("L89").Value = 12833.3335
$ws_BSM.Range("M89").Value = -3607.823
$ws_BSM.Range("N89").Value = -24065.3335
$ws_BSM.Range("H134").Value = 56865.156
$ws_BSM.Range("I134").Value = 59900.89
$ws_BSM.Range("K134").Value = 179702.67
$ws_BSM.Range("M134").Value = -177167.67
$ws_CRP.Range("H58").Value = 24130.5
$ws_CRP.Range("I58").Value = 1840.1
$ws_CRP.Range("K58").Value = 1840.1
$ws_CRP.Range("M58").Value = -1637.1
$ws_CRP.Range("H62").Value = 6168.6665
$ws_CRP.Range("I62").Value = 4500
$ws_CRP.Range("K62").Value = 4500
$ws_CRP.Range("M62").Value = -3876
$ws_CRP.Range("H65").Value = 6168.6665
$ws_CRP.Range("I65").Value = 4500
$ws_CRP.Range("K65").Value = 22500
$ws_CRP.Range("M65").Value = -19380
$ws_CRP.Range("H136").Value = 24130.5
$ws_CRP.Range("I136").Value = 1840.1
$ws_CRP.Range("K136").Value = 5520.299999999999
$ws_CRP.Range("M136").Value = -2970.299999999999
$ws_CUL.Range("H134").Value = 2566
$ws_CUL.Range("I134").Value = 1518.6875
$ws_CUL.Range("J134").Value = 3855
$ws_CUL.Range("K134").Value = 4556.0625
$ws_CUL.Range("L134").Value = 11565
$ws_CUL.Range("M134").Value = 513.9375
$ws_CUL.Range("N134").Value = -21705
$ws_GSM.Range("H80").Value = 3129.6155
$ws_GSM.Range("I80").Value = 3053.7273
$ws_GSM.Range("J80").Value = 3185.2666
$ws_GSM.Range("K80").Value = 3053.7273
$ws_GSM.Range("L80").Value = 3185.2666
$ws_GSM.Range("M80").Value = -2055.7273
$ws_GSM.Range("N80").Value = -5181.2666
$ws_GSM.Range("H83").Value = 3129.6155
$ws_GSM.Range("I83").Value = 3053.7273
$ws_GSM.Range("J83").Value = 3185.2666
$ws_GSM.Range("K83").Value = 15268.6365
$ws_GSM.Range("L83").Value = 15926.333
$ws_GSM.Range("M83").Value = -10276.6365
$ws_GSM.Range("N83").Value = -25910.333
$ws_LTW.Range("H68").Value = 4201.087
$ws_LTW.Range("I68").Value = 2020.6364
$ws_LTW.Range("J68").Value = 6199.8335
$ws_LTW.Range("K68").Value = 2020.6364
$ws_LTW.Range("L68").Value = 6199.8335
$ws_LTW.Range("M68").Value = -1271.6364
$ws_LTW.Range("N68").Value = -7697.8335
$ws_LTW.Range("H71").Value = 4201.087
$ws_LTW.Range("I71").Value = 2020.6364
$ws_LTW.Range("J71").Value = 6199.8335
$ws_LTW.Range("K71").Value = 10103.182
$ws_LTW.Range("L71").Value = 30999.1675
$ws_LTW.Range("M71").Value = -6359.182000000001
$ws_LTW.Range("N71").Value = -38487.1675
$ws_LTW.Range("H82").Value = 3123.3845
$ws_LTW.Range("I82").Value = 3389
$ws_LTW.Range("J82").Value = 2525.75
$ws_LTW.Range("K82").Value = 3389
$ws_LTW.Range("L82").Value = 2525.75
$ws_LTW.Range("M82").Value = -3028
$ws_LTW.Range("N82").Value = -3247.75
$ws_LTW.Range("H85").Value = 3123.3845
$ws_LTW.Range("I85").Value = 3389
$ws_LTW.Range("J85").Value = 2525.75
$ws_LTW.Range("K85").Value = 3389
$ws_LTW.Range("L85").Value = 2525.75
$ws_LTW.Range("M85").Value = -2141
$ws_LTW.Range("N85").Value = -5021.75
$ws_LTW.Range("H132").Value = 2324.7646
$ws_LTW.Range("I132").Value = 1685.8
$ws_LTW.Range("K132").Value = 5057.4
$ws_LTW.Range("M132").Value = -2527.4
$ws_WVR.Range("H27").Value = 32355.6
$ws_WVR.Range("J27").Value = 32355.6
$ws_WVR.Range("L27").Value = 32355.6
$ws_WVR.Range("N27").Value = -32493.6
$ws_WVR.Range("H62").Value = 4249.5
$ws_WVR.Range("I62").Value = 2999
$ws_WVR.Range("K62").Value = 2999
$ws_WVR.Range("M62").Value = -2375
$ws_WVR.Range("H65").Value = 4249.5
$ws_WVR.Range("I65").Value = 2999
$ws_WVR.Range("K65").Value = 14995
$ws_WVR.Range("M65").Value = -11875
$ws_WVR.Range("H108").Value = 32000
$ws_WVR.Range("J108").Value = 32000
$ws_WVR.Range("L108").Value = 32000
$ws_WVR.Range("N108").Value = -39680
